$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.848.27"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "3.842.86"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.24"
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.10"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").Value = "3.842.99"
$ws.Range("E7").Value = "  -1.45%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("E11").Value = "  -1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.78"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "4.480.63"
$ws.Range("E15").Value = "  -1.81%  "

$ws.Range("D16").Value = "3.818.88"
$ws.Range("E16").Value = "  -3.17%  "

$ws.Range("D17").Value = "67.852.95"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("E18").Value = "  +7.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("E20").Value = "  -1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.69"
$ws.Range("E21").Value = "  -3.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.98"
$ws.Range("E22").Value = "  -3.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("E23").Value = "  +1.57%  "

$ws.Range("E24").Value = "  -3.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.06"
$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("E27").Value = "  +1.05%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").Value = "3.987.17"
$ws.Range("E31").Value = "  -1.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.68"
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("E33").Value = "  -3.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.89"
$ws.Range("E34").Value = "  -2.87%  "

$ws.Range("D35").Value = "3.807.75"
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  -2.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.89"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.26"
$ws.Range("E40").Value = "  +9.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("E42").Value = "  -2.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "424.52"
$ws.Range("E43").Value = "  -2.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("E46").Value = "  -2.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.53"
$ws.Range("E47").Value = "  +1.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.34"
$ws.Range("E48").Value = "  +0.87%  "

# Rows 49-51 reordering with updated values
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "39.53"
$ws.Range("E49").Value = "  +1.66%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000266"
$ws.Range("E50").Value = "  +10.66%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0354"
$ws.Range("E51").Value = "  +0.16%  "
